# Milestone 4 further commit
# Updates the AssetList sheet: revises a few "Assets Required" / notes
# cells, fills in newly-added "Status" (Iteration) cells, and moves the
# active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Hurt): reword the required asset, add Status "Iteration 1"
$ws.Range("F6").Value = "Viscous liquid sample, sizzling noise"
$ws.Range("G6").Value = "Iteration 1"

# Row 8 (LandOnGround): reword the required asset, add Status + Notes
$ws.Range("F8").Value = "Splash in mud"
$ws.Range("G8").Value = "Iteration 1"
$ws.Range("H8").Value = "Same multi instrument as Walk, needs something to differentiate from normal walking"

# Row 13 (Menu toggle): reword the required asset, add Status "Iteration 1"
$ws.Range("F13").Value = "Mustard bottle"
$ws.Range("G13").Value = "Iteration 1"

# Row 14 (Button hover): reword the required asset, add Status "Iteration 1"
$ws.Range("F14").Value = "Brief pop"
$ws.Range("G14").Value = "Iteration 1"

# Update the view: scroll so column E is leftmost, and select F14
$win = $excel.ActiveWindow
$win.ScrollColumn = 5
$win.ScrollRow = 1
$ws.Range("F14").Select() | Out-Null
